$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Sheet2 changes
$ws2.Range("B2").Value = 2
$ws2.Range("C2").Value = 100
$ws2.Range("N2").Value = "Yes"

# Sheet1 changes
$ws1.Range("G2").Value = "2:30"
$ws1.Range("G3").Value = "2:30"
$ws1.Range("D3").Value = "D"
$ws1.Range("D3").Select()
